# Add the new "problem 25" building-block rows (126-130) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order chosen to reproduce the shared-string table order of the
# reference commit (strings are appended to the shared-string table the
# first time each distinct value is written).
$ws.Range("B126").Value = "부정방정식 (나)를 풀어 가능한 자연수해를 구합니다."
$ws.Range("A128").Value = "x0008"
$ws.Range("B127").Value = "구해진 자연수 해를 (가)에 대입해서 나머지 자연수들의 관계식을 구합니다."
$ws.Range("C126").Value = "`$\left|a^{2}-b^{2}\right|=5`$; "
$ws.Range("C128").Value = "`$c+d+e=7`$; "
$ws.Range("C127").Value = "`$a+b+c+d+e=12`$; "
$ws.Range("A129").Value = "x0009"
$ws.Range("B128").Value = "방정식을 만족시키는 자연수해의 개수를 중복조합으로 구합니다."
$ws.Range("B129").Value = "방정식을 만족시키는 음이 아닌 정수해의 개수를 중복조합으로 구합니다."
$ws.Range("A130").Value = "x0010"
$ws.Range("B130").Value = "각 경우의 개수를 합해서 문제에서 요구하는 개수를 구합니다. "

$ws.Application.ActiveWindow.ScrollRow = 115
$ws.Range("B131").Select()
